# Auto-generated PowerShell Excel COM-interop script
# Applies the updated vm_pu values per the commit "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.042583576362612
$ws.Range("D2").Value = 1.04621432663889
$ws.Range("E2").Value = 1.050249189732453
$ws.Range("F2").Value = 1.060013728620031
$ws.Range("I2").Value = 1.045043220182748
$ws.Range("J2").Value = 1.047658616253984
$ws.Range("K2").Value = 1.048980389635455
$ws.Range("L2").Value = 1.053003992029427
$ws.Range("M2").Value = 1.062741663782414
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.043482668522081
$ws.Range("D3").Value = 1.046910058106103
$ws.Range("E3").Value = 1.051061187548893
$ws.Range("F3").Value = 1.06093633350908
$ws.Range("I3").Value = 1.045306150462784
$ws.Range("J3").Value = 1.048204456751717
$ws.Range("K3").Value = 1.049488019752257
$ws.Range("L3").Value = 1.053628399485885
$ws.Range("M3").Value = 1.06347833735385
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.044064908164322
$ws.Range("D4").Value = 1.047360543556813
$ws.Range("E4").Value = 1.051587411089929
$ws.Range("F4").Value = 1.061534299556386
$ws.Range("I4").Value = 1.045475204905993
$ws.Range("J4").Value = 1.048557468687856
$ws.Range("K4").Value = 1.049816114064151
$ws.Range("L4").Value = 1.05403258557085
$ws.Range("M4").Value = 1.063955360593248
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.044309792222398
$ws.Range("D5").Value = 1.047549998107094
$ws.Range("E5").Value = 1.051808827043238
$ws.Range("F5").Value = 1.061785917271144
$ws.Range("I5").Value = 1.045546016370767
$ws.Range("J5").Value = 1.048705829920034
$ws.Range("K5").Value = 1.049953953708489
$ws.Range("L5").Value = 1.054202540741627
$ws.Range("M5").Value = 1.064155982676811
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.04435091577716
$ws.Range("D6").Value = 1.047581812447898
$ws.Range("E6").Value = 1.051846014945497
$ws.Range("F6").Value = 1.061828178597434
$ws.Range("I6").Value = 1.045557890717695
$ws.Range("J6").Value = 1.048730737754147
$ws.Range("K6").Value = 1.049977092210679
$ws.Range("L6").Value = 1.054231078999643
$ws.Range("M6").Value = 1.064189672747851
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.044068179884426
$ws.Range("D7").Value = 1.047363074782802
$ws.Range("E7").Value = 1.051590368909899
$ws.Range("F7").Value = 1.061537660771891
$ws.Range("I7").Value = 1.045476152110908
$ws.Range("J7").Value = 1.048559451276027
$ws.Range("K7").Value = 1.04981795624399
$ws.Range("L7").Value = 1.054034856382821
$ws.Range("M7").Value = 1.063958040996761
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.042887331333458
$ws.Range("D8").Value = 1.046449389180143
$ws.Range("E8").Value = 1.050523440247473
$ws.Range("F8").Value = 1.060325323602525
$ws.Range("I8").Value = 1.045132301700033
$ws.Range("J8").Value = 1.047843122918551
$ws.Range("K8").Value = 1.049152022859911
$ws.Range("L8").Value = 1.053214981101254
$ws.Range("M8").Value = 1.06299055356707
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.040810157012517
$ws.Range("D9").Value = 1.044841728025697
$ws.Range("E9").Value = 1.048649621639499
$ws.Range("F9").Value = 1.058196598536248
$ws.Range("I9").Value = 1.044518160277453
$ws.Range("J9").Value = 1.046579504527579
$ws.Range("K9").Value = 1.047975730173234
$ws.Range("L9").Value = 1.051771479471631
$ws.Range("M9").Value = 1.061288428620516
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.039427890323044
$ws.Range("D10").Value = 1.043771636894353
$ws.Range("E10").Value = 1.047404699969162
$ws.Range("F10").Value = 1.056782630626254
$ws.Range("I10").Value = 1.04410323872137
$ws.Range("J10").Value = 1.045736244723695
$ws.Range("K10").Value = 1.047189696206422
$ws.Range("L10").Value = 1.050810039722144
$ws.Range("M10").Value = 1.060155582427223
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.038829965811559
$ws.Range("D11").Value = 1.043308693622598
$ws.Range("E11").Value = 1.046866671976006
$ws.Range("F11").Value = 1.056171615685094
$ws.Range("I11").Value = 1.043922278655166
$ws.Range("J11").Value = 1.045370916319928
$ws.Range("K11").Value = 1.046848912434375
$ws.Range("L11").Value = 1.050393952983269
$ws.Range("M11").Value = 1.059665515935599
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.038607962289305
$ws.Range("D12").Value = 1.043136799469529
$ws.Range("E12").Value = 1.046666980897544
$ws.Range("F12").Value = 1.055944845769729
$ws.Range("I12").Value = 1.043854867942811
$ws.Range("J12").Value = 1.045235189396137
$ws.Range("K12").Value = 1.04672226727047
$ws.Range("L12").Value = 1.050239434545499
$ws.Range("M12").Value = 1.059483554378921
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.038655578603892
$ws.Range("D13").Value = 1.043173668447118
$ws.Range("E13").Value = 1.046709808206372
$ws.Range("F13").Value = 1.055993480142945
$ws.Range("I13").Value = 1.043869336541535
$ws.Range("J13").Value = 1.045264304517972
$ws.Range("K13").Value = 1.046749435916844
$ws.Range("L13").Value = 1.050272577685293
$ws.Range("M13").Value = 1.05952258251796
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.038811613042006
$ws.Range("D14").Value = 1.043294483492355
$ws.Range("E14").Value = 1.046850162248286
$ws.Range("F14").Value = 1.056152866964048
$ws.Range("I14").Value = 1.043916710420093
$ws.Range("J14").Value = 1.045359697648773
$ws.Range("K14").Value = 1.046838445182428
$ws.Range("L14").Value = 1.050381179726135
$ws.Range("M14").Value = 1.059650473481272
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.038907763271138
$ws.Range("D15").Value = 1.043368930100486
$ws.Range("E15").Value = 1.046936659785913
$ws.Range("F15").Value = 1.056251095435427
$ws.Range("I15").Value = 1.043945873331855
$ws.Range("J15").Value = 1.045418468867748
$ws.Range("K15").Value = 1.046893278426584
$ws.Range("L15").Value = 1.050448097662179
$ws.Range("M15").Value = 1.059729280767508
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.039467585476067
$ws.Range("D16").Value = 1.043802369749776
$ws.Range("E16").Value = 1.047440428974294
$ws.Range("F16").Value = 1.05682320798735
$ws.Range("I16").Value = 1.044115221190441
$ws.Range("J16").Value = 1.045760486425632
$ws.Range("K16").Value = 1.047212304022937
$ws.Range("L16").Value = 1.050837658829566
$ws.Range("M16").Value = 1.060188116403634
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.039818910271734
$ws.Range("D17").Value = 1.044074366691559
$ws.Range("E17").Value = 1.047756707267905
$ws.Range("F17").Value = 1.057182412834019
$ws.Range("I17").Value = 1.044221102043436
$ws.Range("J17").Value = 1.045974974522639
$ws.Range("K17").Value = 1.047412307097723
$ws.Range("L17").Value = 1.051082080810931
$ws.Range("M17").Value = 1.060476056882983
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.040023890478557
$ws.Range("D18").Value = 1.044233057677641
$ws.Range("E18").Value = 1.047941286448838
$ws.Range("F18").Value = 1.05739205072127
$ws.Range("I18").Value = 1.044282735435164
$ws.Range("J18").Value = 1.046100063267775
$ws.Range("K18").Value = 1.047528924343621
$ws.Range("L18").Value = 1.051224669454941
$ws.Range("M18").Value = 1.060644052313079
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.040093793300506
$ws.Range("D19").Value = 1.044287173903144
$ws.Range("E19").Value = 1.048004239993497
$ws.Range("F19").Value = 1.057463552091932
$ws.Range("I19").Value = 1.044303729586471
$ws.Range("J19").Value = 1.046142712121568
$ws.Range("K19").Value = 1.047568680795692
$ws.Range("L19").Value = 1.051273292097493
$ws.Range("M19").Value = 1.060701341940938
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.03978121039701
$ws.Range("D20").Value = 1.044045179869206
$ws.Range("E20").Value = 1.047722763307412
$ws.Range("F20").Value = 1.057143861152972
$ws.Range("I20").Value = 1.044209754966723
$ws.Range("J20").Value = 1.045951963889022
$ws.Range("K20").Value = 1.047390852910528
$ws.Range("L20").Value = 1.051055854419472
$ws.Range("M20").Value = 1.06044515896964
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.038765662232395
$ws.Range("D21").Value = 1.04325890470528
$ws.Range("E21").Value = 1.046808827149426
$ws.Range("F21").Value = 1.056105926308918
$ws.Range("I21").Value = 1.043902765345991
$ws.Range("J21").Value = 1.045331607505189
$ws.Range("K21").Value = 1.046812235899769
$ws.Range("L21").Value = 1.050349198164339
$ws.Range("M21").Value = 1.059612810810604
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.038127681205112
$ws.Range("D22").Value = 1.042764910962246
$ws.Range("E22").Value = 1.046235105281014
$ws.Range("F22").Value = 1.055454425348613
$ws.Range("I22").Value = 1.043708626194142
$ws.Range("J22").Value = 1.044941405153422
$ws.Range("K22").Value = 1.046448073051969
$ws.Range("L22").Value = 1.049905096705745
$ws.Range("M22").Value = 1.059089891695782
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03846583595878
$ws.Range("D23").Value = 1.043026750945305
$ws.Range("E23").Value = 1.046539159702276
$ws.Range("F23").Value = 1.055799694503093
$ws.Range("I23").Value = 1.043811649217287
$ws.Range("J23").Value = 1.045148273584068
$ws.Range("K23").Value = 1.046641156748115
$ws.Range("L23").Value = 1.050140503803392
$ws.Range("M23").Value = 1.059367061586836
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.039798245167297
$ws.Range("D24").Value = 1.04405836801518
$ws.Range("E24").Value = 1.047738100815093
$ws.Range("F24").Value = 1.05716128062828
$ws.Range("I24").Value = 1.044214882608919
$ws.Range("J24").Value = 1.045962361461246
$ws.Range("K24").Value = 1.047400547261541
$ws.Range("L24").Value = 1.051067704930761
$ws.Range("M24").Value = 1.06045912026897
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.041346718332811
$ws.Range("D25").Value = 1.045257056903235
$ws.Range("E25").Value = 1.049133299076048
$ws.Range("F25").Value = 1.058746018484095
$ws.Range("I25").Value = 1.044677901947069
$ws.Range("J25").Value = 1.046906334271388
$ws.Range("K25").Value = 1.048280158963487
$ws.Range("L25").Value = 1.052144506836678
$ws.Range("M25").Value = 1.061728139079466
